$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.472.98'
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").Value = '1.746.82'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4456'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3583'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07503'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.95'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.088'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.71'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.007'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.094'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").Value = '1.751.27'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001059'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06398'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.804'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.35%  '
$ws.Range("D23").Value = '27.541.06'
$ws.Range("E23").Value = '  -1.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.087'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").Value = '1.951.24'
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.071'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.12%  '
$ws.Range("E32").Value = '  +4.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09049'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.520'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.90'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02279'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2085'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6342'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05988'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.198'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.384'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.739'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.718'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5859'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.947'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.139'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06841'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.28'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.96%  '
